$d = $word.ActiveDocument

# 1. Remove the table (and its contents: the dice image + the grid of cells).
if ($d.Tables.Count -gt 0) {
    $d.Tables.Item(1).Delete()
}

# 2. Remove the now-unused "Table Grid" style that the table used.
foreach ($styleName in @("Tablaconcuadrcula", "Table Grid")) {
    $removed = $false
    try {
        $s = $d.Styles.Item($styleName)
        if ($s -ne $null) {
            $s.Delete()
            $removed = $true
        }
    } catch {
    }
    if ($removed) { break }
}

# 3. Drop the leading empty paragraph that used to sit right after the table.
$d.Range(0, 1).Delete()

# 4. Remove the two trailing paragraphs entirely (text + paragraph mark).
$d.Content.Find.Execute("Ahora otra.^p", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute("Línea, se ve bien^p", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 2) | Out-Null

# 5. Collapse what's left (the bookmark paragraph + the document's trailing
#    empty paragraph) down to a single paragraph, keeping the bookmark.
$d.Range(0, 2).Delete()

# 6. Insert the restored sentence ("Digame algo Padre") as its own paragraph,
#    spell-check markers and all, right before the bookmark paragraph.
$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Digame</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> algo</w:t></w:r>
<w:r><w:t xml:space="preserve"> Padre</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$d.Range(0, 0).InsertXML($xmlFrag)

# 7. Merge that paragraph into the bookmark paragraph so the bookmark ends up
#    right after "Padre" inside the very same paragraph.
$mergeAt = $d.Content.Text.Length - 2
$d.Range($mergeAt, $mergeAt + 1).Delete()
